$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.077.75"
Set-TextValue "E2" "  -0.47%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.652.43"
Set-TextValue "E3" "  -0.44%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.001"

# Row 5 - BNB
Set-TextValue "D5" "217.40"
Set-TextValue "E5" "  +0.09%  "

# Row 6 - XRP
Set-TextValue "E6" "  +1.79%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.2596"
Set-TextValue "E8" "  -1.75%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.06321"
Set-TextValue "E9" "  +0.67%  "

# Row 10 - Solana
Set-TextValue "E10" "  -2.22%  "

# Row 11 - TRON
Set-TextValue "D11" "0.07793"
Set-TextValue "E11" "  +0.28%  "

# Row 12 - Polkadot
Set-TextValue "D12" "4.515"
Set-TextValue "E12" "  +0.69%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.652.50"
Set-TextValue "E13" "  -0.40%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "1.879.32"
Set-TextValue "E14" "  -0.41%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.5493"
Set-TextValue "E15" "  +0.36%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0₅8204"

# Row 17 - Litecoin
Set-TextValue "D17" "65.58"
Set-TextValue "E17" "  +0.90%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "26.089.45"
Set-TextValue "E18" "  -0.44%  "

# Row 20 - Uniswap
Set-TextValue "D20" "4.583"
Set-TextValue "E20" "  -0.74%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "190.82"
Set-TextValue "E21" "  -0.80%  "

# Row 22 - Avalanche
Set-TextValue "E22" "  -0.19%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.035"
Set-TextValue "E23" "  +0.31%  "

# Row 24 - BinanceUSD
Set-TextValue "E24" "  -0.22%  "

# Row 25 - Monero
Set-TextValue "D25" "144.38"
Set-TextValue "E25" "  +3.51%  "

# Row 26 - Stellar
Set-TextValue "E26" "  +1.31%  "

# Row 27 - Cosmos
Set-TextValue "D27" "7.235"
Set-TextValue "E27" "  -0.78%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "16.04"
Set-TextValue "E28" "  -0.78%  "

# Row 29 - Toncoin
Set-TextValue "D29" "1.431"
Set-TextValue "E29" "  -0.77%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.05835"
Set-TextValue "E30" "  -1.71%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.274"

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "3.551"
Set-TextValue "E32" "  +0.03%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.264"
Set-TextValue "E33" "  -0.62%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.585"
Set-TextValue "E34" "  +0.02%  "

# Row 35 - was MXToken, now HuobiToken
Set-TextValue "B35" "HuobiToken"
Set-TextValue "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "2.411"
Set-TextValue "E35" "  -0.29%  "

# Row 36 - ARBITRUM
Set-TextValue "D36" "0.9455"
Set-TextValue "E36" "  -1.69%  "

# Row 37 - was HuobiToken, now MXToken
Set-TextValue "B37" "MXToken"
Set-TextValue "C37" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D37" "2.780"
Set-TextValue "E37" "  +0.36%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.5740"
Set-TextValue "E38" "  +1.07%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.01605"
Set-TextValue "E39" "  +0.75%  "

# Row 40 - was Quant, now TrustWalletToken
Set-TextValue "B40" "TrustWalletToken"
Set-TextValue "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D40" "0.8447"
Set-TextValue "E40" "  -0.99%  "

# Row 41 - was TrustWalletToken, now Quant
Set-TextValue "B41" "Quant"
Set-TextValue "C41" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D41" "104.32"
Set-TextValue "E41" "  +3.28%  "

# Row 42 - PaxDollar
Set-TextValue "E42" "  -0.11%  "

# Row 43 - FraxShare
Set-TextValue "D43" "5.736"
Set-TextValue "E43" "  -5.06%  "

# Row 44 - Maker
Set-TextValue "D44" "1.030.24"

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.794.73"
Set-TextValue "E45" "  -0.35%  "

# Row 46 - Aave
Set-TextValue "D46" "57.13"
Set-TextValue "E46" "  +1.05%  "

# Row 47 - Frax
Set-TextValue "E47" "  -0.62%  "

# Row 48 - Mantle
Set-TextValue "E48" "  +2.08%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.05141"
Set-TextValue "E49" "  -0.51%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.799"
Set-TextValue "E50" "  -3.04%  "

# Row 51 - RenderToken
Set-TextValue "D51" "1.456"
Set-TextValue "E51" "  +0.21%  "
